# Generate Report for Handback
#
# The a6f5cf38-7a53-43f0-9081-cb5230a1fab6.md file has now been handed
# back in sync with en-US for both locales, so update the status,
# handback timestamps, and clear the stale "version not latest" error.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status columns for the a6f5cf38 row (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: status / handback datetime / error detail ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-22 16:49:40"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).AutoFit()

# --- de-de sheet: status / handback datetime / error detail ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-22 16:49:49"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).AutoFit()
